# Updates the D (Price) and E (Volume 1h) columns in the cryptos sheet
# to reflect the latest GitHub Actions scrape, matching the target diff.
# Values must remain TEXT (inlineStr in the original file) even though many
# look numeric (e.g. "1.003", "225.84") -- Excel's Range.Value setter will
# silently coerce a numeric-looking string into a real number, which would
# both change the cell type and normalize/round the text (e.g. "0.05250"
# would collapse to 0.0525). To prevent that coercion we briefly force the
# cell to Text number-format before writing the value, then restore the
# cell's style to "Normal" so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $addr, $val) {
    $r = $sheet.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" "27.548.53"
Set-TextValue $ws "D3" "1.724.11"
Set-TextValue $ws "E3" "  +4.27%  "
Set-TextValue $ws "D4" "1.003"
Set-TextValue $ws "E4" "  +0.07%  "
Set-TextValue $ws "D5" "225.84"
Set-TextValue $ws "E5" "  +3.42%  "
Set-TextValue $ws "D6" "0.5371"
Set-TextValue $ws "E6" "  +3.28%  "
Set-TextValue $ws "D7" "1.003"
Set-TextValue $ws "E7" "  +0.06%  "
Set-TextValue $ws "D8" "0.2676"
Set-TextValue $ws "E8" "  +1.25%  "
Set-TextValue $ws "D10" "21.75"
Set-TextValue $ws "E10" "  +6.59%  "
Set-TextValue $ws "D11" "0.07719"
Set-TextValue $ws "E11" "  +0.45%  "
Set-TextValue $ws "E12" "  +0.42%  "
Set-TextValue $ws "D13" "1.719.62"
Set-TextValue $ws "E13" "  +5.07%  "
Set-TextValue $ws "D14" "1.961.57"
Set-TextValue $ws "E14" "  +4.42%  "
Set-TextValue $ws "D15" "0.5860"
Set-TextValue $ws "E15" "  +4.76%  "
Set-TextValue $ws "D16" "0.0₅8317"
Set-TextValue $ws "E16" "  +2.04%  "
Set-TextValue $ws "D17" "68.02"
Set-TextValue $ws "E17" "  +3.99%  "
Set-TextValue $ws "D18" "27.564.38"
Set-TextValue $ws "E18" "  +5.60%  "
Set-TextValue $ws "D19" "220.37"
Set-TextValue $ws "E19" "  +14.96%  "
Set-TextValue $ws "E20" "  +0.08%  "
Set-TextValue $ws "D21" "4.734"
Set-TextValue $ws "E21" "  +2.22%  "
Set-TextValue $ws "D22" "10.67"
Set-TextValue $ws "E22" "  +1.85%  "
Set-TextValue $ws "E23" "  +2.90%  "
Set-TextValue $ws "E24" "  +0.05%  "
Set-TextValue $ws "D25" "148.51"
Set-TextValue $ws "E25" "  +3.00%  "
Set-TextValue $ws "D26" "1.747"
Set-TextValue $ws "E26" "  +16.17%  "
Set-TextValue $ws "D27" "0.1237"
Set-TextValue $ws "E27" "  +3.87%  "
Set-TextValue $ws "D28" "7.424"
Set-TextValue $ws "E28" "  +2.87%  "
Set-TextValue $ws "D29" "16.67"
Set-TextValue $ws "E29" "  +4.87%  "
Set-TextValue $ws "D30" "0.05612"
Set-TextValue $ws "E30" "  +2.26%  "
Set-TextValue $ws "D32" "3.558"
Set-TextValue $ws "E32" "  +3.39%  "
Set-TextValue $ws "D33" "3.459"
Set-TextValue $ws "E33" "  +2.93%  "
Set-TextValue $ws "E34" "  +6.80%  "
Set-TextValue $ws "E35" "  +1.81%  "
Set-TextValue $ws "D36" "0.9612"
Set-TextValue $ws "E36" "  +1.32%  "
Set-TextValue $ws "D37" "2.427"
Set-TextValue $ws "E37" "  +0.19%  "
Set-TextValue $ws "D38" "0.5966"
Set-TextValue $ws "E38" "  +5.81%  "
Set-TextValue $ws "D39" "0.01653"
Set-TextValue $ws "E39" "  +4.55%  "
Set-TextValue $ws "D40" "5.931"
Set-TextValue $ws "E40" "  +1.41%  "
Set-TextValue $ws "D41" "0.8568"
Set-TextValue $ws "E41" "  +3.58%  "
Set-TextValue $ws "D42" "1.056.24"
Set-TextValue $ws "E42" "  +2.72%  "
Set-TextValue $ws "E43" "  +0.07%  "
Set-TextValue $ws "D44" "101.86"
Set-TextValue $ws "E44" "  +0.57%  "
Set-TextValue $ws "D45" "1.868.15"
Set-TextValue $ws "E45" "  +4.14%  "
Set-TextValue $ws "E46" "  +6.03%  "
Set-TextValue $ws "D47" "59.13"
Set-TextValue $ws "E47" "  +2.74%  "
Set-TextValue $ws "D48" "8.206"
Set-TextValue $ws "E48" "  +2.92%  "
Set-TextValue $ws "D49" "0.4434"
Set-TextValue $ws "E49" "  +2.27%  "
Set-TextValue $ws "E50" "  +0.24%  "
Set-TextValue $ws "D51" "0.05250"
Set-TextValue $ws "E51" "  +1.63%  "
